# Add two new Protocol test-case rows (protocol_031, protocol_032) covering
# Chinese-character / encoding variable support, appended after the existing
# last row (31) of Sheet1.
#
# New rows are created via Rows.Insert() (copying the immediately preceding
# row's cell formatting down, the way Excel normally extends a table) rather
# than writing into previously-blank rows, so the inherited style indexes for
# the untouched columns (A/B/D/I/J/K) line up with the neighbouring rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 32: protocol_031 ----------------------------------------------
$ws.Rows.Item(32).Insert()
$ws.Cells.Item(32, 5).Clear()   # Table_used  - not used by this case
$ws.Cells.Item(32, 8).Clear()   # Op_sql      - not used by this case

$ws.Cells.Item(32, 1).Value = "protocol_031"
$ws.Cells.Item(32, 2).Value = "y"
$ws.Cells.Item(32, 3).Value = "查看编码变量"
$ws.Cells.Item(32, 3).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "Protocol"
$ws.Cells.Item(32, 9).Value = "show variables like '%character%'"
$ws.Cells.Item(32, 10).Value = "src/test/resources/io.dingodb.test/testdata/mysqlcases/protocol/expectedresult/protocol_031.csv"
$ws.Cells.Item(32, 11).Value = "csv_containsAll"

# ---- Row 33: protocol_032 ----------------------------------------------
$ws.Rows.Item(33).Insert()

$ws.Cells.Item(33, 1).Value = "protocol_032"
$ws.Cells.Item(33, 2).Value = "y"
$ws.Cells.Item(33, 3).Value = "设置编码变量"
$ws.Cells.Item(33, 3).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "Protocol"
$ws.Cells.Item(33, 8).Value = "set names utf8"
$ws.Cells.Item(33, 9).Value = "show variables like '%character%'"
$ws.Cells.Item(33, 10).Value = "src/test/resources/io.dingodb.test/testdata/mysqlcases/protocol/expectedresult/protocol_032.csv"
$ws.Cells.Item(33, 11).Value = "csv_containsAll"

# ---- Scroll / selection to match the authored view ----------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 25
$win.ScrollColumn = 1
$ws.Range("F40").Select() | Out-Null
